$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting the existing columns
# (A->B, B->C, C->D, D->E) one place to the right.
$ws.Columns("A").Insert()

# The new header cell (B1) should look like the other header cells
# (bold, thin box border, centered/top aligned) - copy that formatting
# from the neighboring header cell instead of rebuilding it by hand so
# the existing style is reused rather than a new one being minted.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = "segments"

# The new column A (segment index) should look like the adjacent label
# column (also bold / bordered / centered-top) - copy that formatting too.
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill column A with the 0-based segment index for each data row.
for ($i = 0; $i -lt 19; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
